$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name stored in workbook.xml <sheet name="...">)
$ws.Name = "alpha3F"

# Update the slightly-changed floating point values in row 13
$ws.Range("D13").Value = 0.9924990944342341
$ws.Range("J13").Value = 0.9924990944342341
$ws.Range("K13").Value = 0.9923748021635791
